$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (was Densidade, now Furo)
$ws.Range("A2").Value = "Furo"
$ws.Range("B2").Value = 50
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 59.79981992270027
$ws.Range("E2").Value = 40.20018007729973
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = 0.95
$ws.Range("H2").Value = "Continuous"

# Add new row 3 (Pino)
$ws.Range("A3").Value = "Pino"
$ws.Range("B3").Value = 40
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 45.87989195362016
$ws.Range("E3").Value = 34.12010804637984
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = 0.95
$ws.Range("H3").Value = "Continuous"
